$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as plain text, preserving its default (Normal) style
# so that numeric-looking strings (e.g. prices, percentages) are not
# auto-converted to numbers by Excel's type inference.
function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = "Normal"
}

Set-TextValue "D2" "293.96"
Set-TextValue "E2" "1.39%"
Set-TextValue "D3" "31.17"
Set-TextValue "E3" "1.46%"
Set-TextValue "D4" "4.933"
Set-TextValue "E4" "1.03%"
Set-TextValue "D5" "0.07328"
Set-TextValue "E5" "1.77%"
Set-TextValue "D6" "2.278"
Set-TextValue "E6" "28.09%"
Set-TextValue "D7" "7.745"
Set-TextValue "E7" "0.91%"
Set-TextValue "E8" "0.32%"
Set-TextValue "D9" "0.9089"
Set-TextValue "E9" "1.45%"
Set-TextValue "D10" "0.1692"
Set-TextValue "E10" "2.44%"
Set-TextValue "D11" "0.08069"
Set-TextValue "E11" "8.31%"
Set-TextValue "D12" "0.08076"
Set-TextValue "E12" "0.26%"
Set-TextValue "D13" "0.03100"
Set-TextValue "E13" "3.68%"
Set-TextValue "E14" "0.83%"
Set-TextValue "D15" "0.001519"
Set-TextValue "E15" "1.51%"
Set-TextValue "D16" "0.005712"
Set-TextValue "E16" "0.56%"
Set-TextValue "D17" "3.476"
Set-TextValue "E17" "0.49%"
Set-TextValue "D18" "2.078"
Set-TextValue "E18" "-1.05%"
Set-TextValue "E19" "1.51%"
Set-TextValue "E20" "0.39%"
Set-TextValue "D21" "3.975"
Set-TextValue "E21" "-9.57%"
Set-TextValue "D23" "0.04546"
Set-TextValue "E23" "1.51%"
Set-TextValue "D24" "0.001212"
Set-TextValue "E24" "0.01%"
Set-TextValue "D25" "0.004494"
Set-TextValue "E25" "11.73%"
Set-TextValue "D26" "0.0001303"
Set-TextValue "E26" "4.13%"
Set-TextValue "D39" "0.01601"
Set-TextValue "E39" "-2.72%"
Set-TextValue "D40" "0.04444"
Set-TextValue "E40" "2.42%"
Set-TextValue "D41" "0.007364"
Set-TextValue "E41" "-0.23%"
Set-TextValue "D42" "0.1330"
Set-TextValue "E42" "1.48%"
Set-TextValue "D43" "0.008630"
Set-TextValue "D44" "0.001970"
Set-TextValue "E44" "-1.90%"
Set-TextValue "D45" "0.009527"
Set-TextValue "E45" "-6.41%"
Set-TextValue "D46" "0.00005980"
Set-TextValue "E46" "2.94%"
Set-TextValue "E47" "0.01%"
Set-TextValue "E48" "2.17%"
Set-TextValue "E49" "-3.46%"
Set-TextValue "E50" "0.01%"
Set-TextValue "E51" "0.01%"
